$d = $word.ActiveDocument

# Locate the paragraph containing "Michael Gallagher" and insert a new
# paragraph right after it containing "Qilin Yang". The existing
# "_GoBack" bookmark that trails "Michael Gallagher" will naturally end
# up following the new paragraph's text because InsertParagraphAfter
# splits at the end of the found range, before the bookmark markers.

$range = $d.Content
$range.Find.Execute("Michael Gallagher") | Out-Null

$range.Collapse(0)
$range.InsertParagraphAfter()
$range.Collapse(0)
$range.InsertAfter("Qilin Yang")

$d.Save()
